# Ata de reunião - update responsible-person names on the "Planilha1" sheet.
# "Tânia" -> "Colaborador F" and "Flávio" -> "Diretor A" in column C,
# then leave the selection on C12 (matches the author's last selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Rows where column C currently reads "Tânia"
$ws.Range("C8").Value  = "Colaborador F"
$ws.Range("C9").Value  = "Colaborador F"
$ws.Range("C12").Value = "Colaborador F"
$ws.Range("C13").Value = "Colaborador F"
$ws.Range("C14").Value = "Colaborador F"
$ws.Range("C15").Value = "Colaborador F"
$ws.Range("C16").Value = "Colaborador F"
$ws.Range("C17").Value = "Colaborador F"
$ws.Range("C18").Value = "Colaborador F"

# Row where column C currently reads "Flávio"
$ws.Range("C10").Value = "Diretor A"

# Update the active selection to match the saved view state
$ws.Range("C12").Select()
